$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.523724666666666
$ws.Range("H2").Value = 13.571174
$ws.Range("I2").Value = 0.4806607624766543
$ws.Range("J2").Value = 0.4806607624766543
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 91.67302333333333
$ws.Range("N2").Value = 275.01907
$ws.Range("O2").Value = 0.8966431814716052
$ws.Range("P2").Value = 0.896643181471605
$ws.Range("Q2").Value = 414.7035169209088
$ws.Range("R2").Value = 3732.33165228818
$ws.Range("S2").Value = 0.4309811952756349
$ws.Range("T2").Value = 0.4309811952756348

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.523724666666666
$ws.Range("H3").Value = 13.571174
$ws.Range("I3").Value = 0.4806607624766543
$ws.Range("J3").Value = 0.4806607624766543
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.548386
$ws.Range("N3").Value = 25.645158
$ws.Range("O3").Value = 0.08361076945850333
$ws.Range("P3").Value = 0.08361076945850332
$ws.Range("Q3").Value = 38.670544608388
$ws.Range("R3").Value = 348.034901475492
$ws.Range("S3").Value = 0.04018841619918398
$ws.Range("T3").Value = 0.04018841619918397

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.523724666666666
$ws.Range("H4").Value = 13.571174
$ws.Range("I4").Value = 0.4806607624766543
$ws.Range("J4").Value = 0.4806607624766543
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.018841
$ws.Range("N4").Value = 6.056523
$ws.Range("O4").Value = 0.01974604906989159
$ws.Range("P4").Value = 0.01974604906989159
$ws.Range("Q4").Value = 9.132680829778
$ws.Range("R4").Value = 82.194127468002
$ws.Range("S4").Value = 0.009491151001835525
$ws.Range("T4").Value = 0.009491151001835523

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.109481
$ws.Range("H5").Value = 6.328443
$ws.Range("I5").Value = 0.2241393587371326
$ws.Range("J5").Value = 0.2241393587371326
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 91.67302333333333
$ws.Range("N5").Value = 275.01907
$ws.Range("O5").Value = 0.8966431814716052
$ws.Range("P5").Value = 0.896643181471605
$ws.Range("Q5").Value = 193.3825009342233
$ws.Range("R5").Value = 1740.44250840801
$ws.Range("S5").Value = 0.200973027711068
$ws.Range("T5").Value = 0.2009730277110679

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.109481
$ws.Range("H6").Value = 6.328443
$ws.Range("I6").Value = 0.2241393587371326
$ws.Range("J6").Value = 0.2241393587371326
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.548386
$ws.Range("N6").Value = 25.645158
$ws.Range("O6").Value = 0.08361076945850333
$ws.Range("P6").Value = 0.08361076945850332
$ws.Range("Q6").Value = 18.032657847666
$ws.Range("R6").Value = 162.293920628994
$ws.Range("S6").Value = 0.01874046424994717
$ws.Range("T6").Value = 0.01874046424994716

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.109481
$ws.Range("H7").Value = 6.328443
$ws.Range("I7").Value = 0.2241393587371326
$ws.Range("J7").Value = 0.2241393587371326
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.018841
$ws.Range("N7").Value = 6.056523
$ws.Range("O7").Value = 0.01974604906989159
$ws.Range("P7").Value = 0.01974604906989159
$ws.Range("Q7").Value = 4.258706731521
$ws.Range("R7").Value = 38.328360583689
$ws.Range("S7").Value = 0.004425866776117455
$ws.Range("T7").Value = 0.004425866776117454

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.778265
$ws.Range("H8").Value = 8.334795
$ws.Range("I8").Value = 0.2951998787862131
$ws.Range("J8").Value = 0.2951998787862131
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 91.67302333333333
$ws.Range("N8").Value = 275.01907
$ws.Range("O8").Value = 0.8966431814716052
$ws.Range("P8").Value = 0.896643181471605
$ws.Range("Q8").Value = 254.6919521711833
$ws.Range("R8").Value = 2292.22756954065
$ws.Range("S8").Value = 0.2646889584849023
$ws.Range("T8").Value = 0.2646889584849023

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.778265
$ws.Range("H9").Value = 8.334795
$ws.Range("I9").Value = 0.2951998787862131
$ws.Range("J9").Value = 0.2951998787862131
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.548386
$ws.Range("N9").Value = 25.645158
$ws.Range("O9").Value = 0.08361076945850333
$ws.Range("P9").Value = 0.08361076945850332
$ws.Range("Q9").Value = 23.74968163029
$ws.Range("R9").Value = 213.74713467261
$ws.Range("S9").Value = 0.02468188900937219
$ws.Range("T9").Value = 0.02468188900937219

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.778265
$ws.Range("H10").Value = 8.334795
$ws.Range("I10").Value = 0.2951998787862131
$ws.Range("J10").Value = 0.2951998787862131
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.018841
$ws.Range("N10").Value = 6.056523
$ws.Range("O10").Value = 0.01974604906989159
$ws.Range("P10").Value = 0.01974604906989159
$ws.Range("Q10").Value = 5.608875290865
$ws.Range("R10").Value = 50.479877617785
$ws.Range("S10").Value = 0.005829031291938614
$ws.Range("T10").Value = 0.005829031291938613
